$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 32
$ws.Range("H32").Value = 17651138
$ws.Range("I32").Value = 25006598
$ws.Range("J32").Value = 11112950
$ws.Range("K32").Value = 25006598
$ws.Range("L32").Value = 11112950
$ws.Range("M32").Value = -25006272
$ws.Range("N32").Value = -11113602
# Row 116
$ws.Range("H116").Value = 5171.278
$ws.Range("J116").Value = 4820.2144
$ws.Range("L116").Value = 4820.2144
$ws.Range("N116").Value = -11704.2144
# Row 123
$ws.Range("H123").Value = 200000
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 200000
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 200000
$ws.Range("M123").ClearContents()
$ws.Range("N123").Value = -209800
# Row 132
$ws.Range("H132").Value = 5346.672
$ws.Range("I132").Value = 2899.2744
$ws.Range("J132").Value = 14948
$ws.Range("K132").Value = 8697.823199999999
$ws.Range("L132").Value = 44844
$ws.Range("M132").Value = -6167.823199999999
$ws.Range("N132").Value = -49904
# Row 133
$ws.Range("H133").Value = 100000
$ws.Range("J133").Value = 100000
$ws.Range("L133").Value = 100000
$ws.Range("N133").Value = -110120
# Row 134
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
# Row 135
$ws.Range("H135").Value = 1544.1786
$ws.Range("I135").Value = 1057.1818
$ws.Range("J135").Value = 3329.8333
$ws.Range("K135").Value = 9514.636200000001
$ws.Range("L135").Value = 29968.4997
$ws.Range("M135").Value = -6979.636200000001
$ws.Range("N135").Value = -35038.4997
# Row 136
$ws.Range("H136").Value = 132871.5
$ws.Range("J136").Value = 132871.5
$ws.Range("L136").Value = 132871.5
$ws.Range("N136").Value = -143071.5
# Row 139
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
# Row 140
$ws.Range("H140").Value = 100918
$ws.Range("J140").Value = 100918
$ws.Range("L140").Value = 100918
$ws.Range("N140").Value = -111278

$ws = $wb.Worksheets.Item("ARM")
# Row 29
$ws.Range("H29").Value = 23598
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 23598
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 23598
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -24214
# Row 32
$ws.Range("H32").Value = 3934.8247
$ws.Range("I32").Value = 3438.043
$ws.Range("K32").Value = 3438.043
$ws.Range("M32").Value = -3151.043
# Row 70
$ws.Range("H70").Value = 49999
$ws.Range("J70").Value = 49999
$ws.Range("L70").Value = 49999
$ws.Range("N70").Value = -50539
# Row 73
$ws.Range("H73").Value = 49999
$ws.Range("J73").Value = 49999
$ws.Range("L73").Value = 49999
$ws.Range("N73").Value = -51871
# Row 80
$ws.Range("H80").Value = 48498.832
$ws.Range("J80").Value = 48498.832
$ws.Range("L80").Value = 48498.832
$ws.Range("N80").Value = -50494.832
# Row 83
$ws.Range("H83").Value = 48498.832
$ws.Range("J83").Value = 48498.832
$ws.Range("L83").Value = 145496.496
$ws.Range("N83").Value = -155480.496
# Row 101
$ws.Range("H101").Value = 28799.4
$ws.Range("J101").Value = 28799.4
$ws.Range("L101").Value = 28799.4
$ws.Range("N101").Value = -35289.4
# Row 119
$ws.Range("H119").Value = 28249.5
$ws.Range("J119").Value = 28249.5
$ws.Range("L119").Value = 28249.5
$ws.Range("N119").Value = -37925.5
# Row 139
$ws.Range("H139").Value = 114983
$ws.Range("J139").Value = 114983
$ws.Range("L139").Value = 114983
$ws.Range("N139").Value = -125263
# Row 140
$ws.Range("H140").Value = 111864.875
$ws.Range("J140").Value = 111864.875
$ws.Range("L140").Value = 111864.875
$ws.Range("N140").Value = -122224.875
# Row 141
$ws.Range("H141").Value = 136075.83
$ws.Range("J141").Value = 136075.83
$ws.Range("L141").Value = 136075.83
$ws.Range("N141").Value = -146435.83

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 1790.4348
$ws.Range("I86").Value = 1455
$ws.Range("J86").Value = 2998
$ws.Range("K86").Value = 1455
$ws.Range("L86").Value = 2998
$ws.Range("M86").Value = -332
$ws.Range("N86").Value = -5244
# Row 89
$ws.Range("H89").Value = 1790.4348
$ws.Range("I89").Value = 1455
$ws.Range("J89").Value = 2998
$ws.Range("K89").Value = 7275
$ws.Range("L89").Value = 14990
$ws.Range("M89").Value = -1659
$ws.Range("N89").Value = -26222
# Row 132
$ws.Range("H132").Value = 133654.1
$ws.Range("J132").Value = 133654.1
$ws.Range("L132").Value = 133654.1
$ws.Range("N132").Value = -143774.1
# Row 135
$ws.Range("H135").Value = 49999
$ws.Range("J135").Value = 49999
$ws.Range("L135").Value = 49999
$ws.Range("N135").Value = -60139
# Row 137
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
# Row 138
$ws.Range("H138").Value = 123442.29
$ws.Range("J138").Value = 123442.29
$ws.Range("L138").Value = 123442.29
$ws.Range("N138").Value = -133722.29
# Row 140
$ws.Range("H140").Value = 130300
$ws.Range("J140").Value = 130300
$ws.Range("L140").Value = 130300
$ws.Range("N140").Value = -140660

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 66.588234
$ws.Range("I7").Value = 41.23077
$ws.Range("J7").Value = 149
$ws.Range("K7").Value = 41.23077
$ws.Range("L7").Value = 149
$ws.Range("M7").Value = 71.76922999999999
$ws.Range("N7").Value = -375
# Row 31
$ws.Range("H31").Value = 1568.0869
$ws.Range("I31").Value = 1568.0869
$ws.Range("K31").Value = 1568.0869
$ws.Range("M31").Value = -1273.0869
# Row 34
$ws.Range("H34").Value = 1568.0869
$ws.Range("I34").Value = 1568.0869
$ws.Range("K34").Value = 1568.0869
$ws.Range("M34").Value = -1366.0869
# Row 134
$ws.Range("H134").Value = 4010.1177
$ws.Range("I134").Value = 4139.452
$ws.Range("J134").Value = 3406.5557
$ws.Range("K134").Value = 12418.356
$ws.Range("L134").Value = 10219.6671
$ws.Range("M134").Value = -9883.356
$ws.Range("N134").Value = -15289.6671
# Row 135
$ws.Range("H135").Value = 169281.78
$ws.Range("J135").Value = 169281.78
$ws.Range("L135").Value = 169281.78
$ws.Range("N135").Value = -179421.78
# Row 138
$ws.Range("H138").Value = 182712.5
$ws.Range("J138").Value = 182712.5
$ws.Range("L138").Value = 182712.5
$ws.Range("N138").Value = -192992.5
# Row 140
$ws.Range("H140").Value = 114676.62
$ws.Range("J140").Value = 118488.3
$ws.Range("L140").Value = 118488.3
$ws.Range("N140").Value = -128848.3

$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 2103266.2
$ws.Range("I131").Value = 7354689
$ws.Range("J131").Value = 2697.2
$ws.Range("K131").Value = 22064067
$ws.Range("L131").Value = 8091.599999999999
$ws.Range("M131").Value = -22059027
$ws.Range("N131").Value = -18171.6
# Row 139
$ws.Range("H139").Value = 3681.55
$ws.Range("I139").Value = 1678.0667
$ws.Range("K139").Value = 5034.2001
$ws.Range("M139").Value = 105.7999
# Row 141
$ws.Range("H141").Value = 15092.125
$ws.Range("I141").Value = 15092.125
$ws.Range("K141").Value = 45276.375
$ws.Range("M141").Value = -40096.375

$ws = $wb.Worksheets.Item("GSM")
# Row 133
$ws.Range("H133").Value = 99999
$ws.Range("J133").Value = 99999
$ws.Range("L133").Value = 99999
$ws.Range("N133").Value = -110119
# Row 135
$ws.Range("H135").Value = 130973.5
$ws.Range("J135").Value = 130973.5
$ws.Range("L135").Value = 130973.5
$ws.Range("N135").Value = -141113.5
# Row 140
$ws.Range("H140").Value = 71943
$ws.Range("I140").Value = 20000
$ws.Range("J140").Value = 82331.60000000001
$ws.Range("K140").Value = 20000
$ws.Range("L140").Value = 82331.60000000001
$ws.Range("M140").Value = -14820
$ws.Range("N140").Value = -92691.60000000001
# Row 141
$ws.Range("H141").Value = 52136.715
$ws.Range("I141").Value = 30390
$ws.Range("J141").Value = 55761.168
$ws.Range("K141").Value = 30390
$ws.Range("L141").Value = 55761.168
$ws.Range("M141").Value = -25210
$ws.Range("N141").Value = -66121.16800000001

$ws = $wb.Worksheets.Item("LTW")
# Row 98
$ws.Range("H98").Value = 95499.25
$ws.Range("J98").Value = 62498.5
$ws.Range("L98").Value = 62498.5
$ws.Range("N98").Value = -68488.5
# Row 132
$ws.Range("H132").Value = 72981.17
$ws.Range("I132").Value = 87574.05499999999
$ws.Range("J132").Value = 3665
$ws.Range("K132").Value = 262722.165
$ws.Range("L132").Value = 10995
$ws.Range("M132").Value = -260192.165
$ws.Range("N132").Value = -16055
# Row 136
$ws.Range("H136").Value = 4564.396
$ws.Range("I136").Value = 3457.738
$ws.Range("J136").Value = 12311
$ws.Range("K136").Value = 10373.214
$ws.Range("L136").Value = 36933
$ws.Range("M136").Value = -7823.214
$ws.Range("N136").Value = -42033

$ws = $wb.Worksheets.Item("WVR")
# Row 46
$ws.Range("H46").Value = 80213.664
$ws.Range("J46").Value = 80213.664
$ws.Range("L46").Value = 80213.664
$ws.Range("N46").Value = -80675.664
# Row 126
$ws.Range("H126").Value = 3676.6843
$ws.Range("I126").Value = 3621.6667
$ws.Range("J126").Value = 3771
$ws.Range("K126").Value = 10865.0001
$ws.Range("L126").Value = 11313
$ws.Range("M126").Value = -8395.000100000001
$ws.Range("N126").Value = -16253
# Row 133
$ws.Range("H133").Value = 96100
$ws.Range("J133").Value = 96100
$ws.Range("L133").Value = 96100
$ws.Range("N133").Value = -106220
# Row 134
$ws.Range("H134").Value = 80213.664
$ws.Range("J134").Value = 80213.664
$ws.Range("L134").Value = 240640.992
$ws.Range("N134").Value = -245710.992
# Row 136
$ws.Range("H136").Value = 1519.2609
$ws.Range("I136").Value = 1583.8813
$ws.Range("K136").Value = 4751.6439
$ws.Range("M136").Value = -2201.6439
# Row 140
$ws.Range("H140").Value = 43964.25
$ws.Range("J140").Value = 43964.25
$ws.Range("L140").Value = 43964.25
$ws.Range("N140").Value = -54324.25

